$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("R9").Value = 1.870901181057595
$ws.Range("S9").Value = 0.4580182439559286
$ws.Range("R12").Value = 57.58568459491445
$ws.Range("S12").Value = 1.767791509175992
$ws.Range("R13").Value = 1.047927592285356
$ws.Range("S13").Value = 0.3113145974078815
$ws.Range("R14").Value = 8.632235739018183
$ws.Range("S14").Value = 0.9837271029525712
$ws.Range("R15").Value = 54.77386105131029
$ws.Range("S15").Value = 1.746430710391437
$ws.Range("R16").Value = 473.1462862645127
$ws.Range("S16").Value = 2.675912353311067
$ws.Range("R17").Value = 976.0242192382194
$ws.Range("S17").Value = 2.989905329482484
$ws.Range("R19").Value = 46.73820180873648
$ws.Range("S19").Value = 1.678866056077288
$ws.Range("R20").Value = 316.1329298790346
$ws.Range("S20").Value = 2.501241339858165
$ws.Range("R21").Value = 1.441722400822429
$ws.Range("S21").Value = 0.3876962875187916
$ws.Range("R22").Value = 3.625693341108359
$ws.Range("S22").Value = 0.6651768379722269
$ws.Range("R23").Value = 1.63192857169589
$ws.Range("S23").Value = 0.4202740987187935
$ws.Range("R24").Value = 130.120275196242
$ws.Range("S24").Value = 2.117669852063998
$ws.Range("R25").Value = 43.931163323227
$ws.Range("S25").Value = 1.652547663158297
$ws.Range("R26").Value = 79.66017359891231
$ws.Range("S26").Value = 1.906659152372425
$ws.Range("R28").Value = 3.150393046867427
$ws.Range("S28").Value = 0.6180892268332367
$ws.Range("R29").Value = 28.158123859662
$ws.Range("S29").Value = 1.464759576449468
$ws.Range("R30").Value = 0.9250166959981613
$ws.Range("S30").Value = 0.2844345005709749
$ws.Range("R31").Value = 1.89488508297257
$ws.Range("S31").Value = 0.4616313284014649
$ws.Range("R32").Value = 1.658182399630851
$ws.Range("S32").Value = 0.4245847781248915
$ws.Range("R33").Value = 2.423023938719238
$ws.Range("S33").Value = 0.5344099363617495
$ws.Range("R34").Value = 0.5407268426004352
$ws.Range("S34").Value = 0.1877256489281545
$ws.Range("R35").Value = 1.26630512812856
$ws.Range("S35").Value = 0.3553183814944465
$ws.Range("R36").Value = 3.214746706901876
$ws.Range("S36").Value = 0.6247714800046056
$ws.Range("R37").Value = 3.885215271104165
$ws.Range("S37").Value = 0.6888837060266939
$ws.Range("R38").Value = 1.767192603285276
$ws.Range("S38").Value = 0.4420393881373103
$ws.Range("R39").Value = 0.886071812742895
$ws.Range("S39").Value = 0.2755582246071517
$ws.Range("R40").Value = 9.524342519382833
$ws.Range("S40").Value = 1.022194973949702
$ws.Range("R42").Value = 2.0000993536122
$ws.Range("S42").Value = 0.4771356373900139
$ws.Range("R43").Value = 0.8247231992570295
$ws.Range("S43").Value = 0.2611969936403649
$ws.Range("R44").Value = 3.714931372753841
$ws.Range("S44").Value = 0.6734753758323673
$ws.Range("R48").Value = 0.5620358249665113
$ws.Range("S48").Value = 0.1936909901090488
$ws.Range("R55").Value = 0.5184274680391626
$ws.Range("S55").Value = 0.1813940514517223
$ws.Range("R57").Value = 1.252623087665084
$ws.Range("S57").Value = 0.3526885309866088
$ws.Range("R58").Value = 3.165109855994819
$ws.Range("S58").Value = 0.6196264605382849
$ws.Range("R62").Value = 0.3253443451059597
$ws.Range("S62").Value = 0.1223287294051758
$ws.Range("R76").Value = 0.7092071748480213
$ws.Range("S76").Value = 0.2327947072146441
$ws.Range("R85").Value = 8.680893291705946
$ws.Range("S85").Value = 0.9859154331099647
$ws.Range("R104").Value = 11.70653431521832
$ws.Range("S104").Value = 1.104027113652473
$ws.Range("R107").Value = 0.7106895364787947
$ws.Range("S107").Value = 0.2331711990043615
$ws.Range("R112").Value = 1.389294654818555
$ws.Range("S112").Value = 0.3782697115208322
$ws.Range("R117").Value = 0.4087187001793404
$ws.Range("S117").Value = 0.1488242797264549
$ws.Range("R118").Value = 0.5582092593561443
$ws.Range("S118").Value = 0.1926257807312954
$ws.Range("R121").Value = 0.9355995287148497
$ws.Range("S121").Value = 0.2868155076933723
$ws.Range("R122").Value = 46.59038322156842
$ws.Range("S122").Value = 1.677519201969222
$ws.Range("R123").Value = 1.185641823985728
$ws.Range("S123").Value = 0.3395789926495084
$ws.Range("R124").Value = 7.439323364126635
$ws.Range("S124").Value = 0.9263076277860751
$ws.Range("R127").Value = 2.882365397774951
$ws.Range("S127").Value = 0.5890964076029308
$ws.Range("R128").Value = 3.478481813310682
$ws.Range("S128").Value = 0.6511308149254167
$ws.Range("R129").Value = 2.957372605549809
$ws.Range("S129").Value = 0.5974069430928867
$ws.Range("R130").Value = 22.09593464075914
$ws.Range("S130").Value = 1.363535541863407
$ws.Range("R131").Value = 16.58962251844254
$ws.Range("S131").Value = 1.245256519393987
$ws.Range("R133").Value = 0.3701862945602766
$ws.Range("S133").Value = 0.1367796191265753
$ws.Range("R138").Value = 1.780848398445973
$ws.Range("S138").Value = 0.4441773133830383
$ws.Range("R139").Value = 2.442220573432152
$ws.Range("S139").Value = 0.5368386959610926
$ws.Range("R140").Value = 2.640682812932694
$ws.Range("S140").Value = 0.5611828435608059
$ws.Range("R143").Value = 47.72146225222763
$ws.Range("S143").Value = 1.687720314077689
$ws.Range("R147").Value = 1.097193358193166
$ws.Range("S147").Value = 0.3216384736341777
$ws.Range("R148").Value = 2.295168426947501
$ws.Range("S148").Value = 0.5178776177234982
$ws.Range("R149").Value = 35.26361622431453
$ws.Range("S149").Value = 1.559471109962378
$ws.Range("R150").Value = 3.299184985570151
$ws.Range("S150").Value = 0.6333861323661171
$ws.Range("R153").Value = 1.305660762531532
$ws.Range("S153").Value = 0.3627954088569834
$ws.Range("R154").Value = 3.597280978501602
$ws.Range("S154").Value = 0.6625010479494922
$ws.Range("R155").Value = 15.87296903244873
$ws.Range("S155").Value = 1.227191509438888
$ws.Range("R158").Value = 0.4658554472824853
$ws.Range("S158").Value = 0.1660911452417727
$ws.Range("R162").Value = 1.712729579514816
$ws.Range("S162").Value = 0.4334065029477489
$ws.Range("R166").Value = 0.8683467697500011
$ws.Range("S166").Value = 0.2714574854980195
$ws.Range("R168").Value = 2.373991978051463
$ws.Range("S168").Value = 0.5281440456828261
$ws.Range("R170").Value = 0.7423734953297659
$ws.Range("S170").Value = 0.2411412560528974
$ws.Range("R172").Value = 19.2664548137769
$ws.Range("S172").Value = 1.306777784732914
